$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (old D:K shifts to F:M) across the sheet.
$ws.Columns("D:E").Insert()

# The newly inserted columns come in with default formatting; copy the
# number formatting from column F (which used to be column D) onto the
# two new columns so dates/numbers keep the right style per row.
$ws.Range("F7:F102").Copy()
$ws.Range("D7:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the two new columns (D = latest quarter, E = prior quarter)
# with the newly reported figures for each line item.

# Income Statement - Period Ending (row 7)
$ws.Range("D7").Value = 43466
$ws.Range("E7").Value = 43375

# Total Revenue
$ws.Range("D8").Value = 113200
$ws.Range("E8").Value = 116700

# Cost of Revenue
$ws.Range("D9").Value = 83100
$ws.Range("E9").Value = 84600

# Gross Profit
$ws.Range("D10").Value = 30100
$ws.Range("E10").Value = 32100

# Research Development (row 12) stays "NA"
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"

# Selling General and Administrative
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0

# Non Recurring
$ws.Range("D14").Value = 800
$ws.Range("E14").Value = 1800

# Others
$ws.Range("D15").Value = 5500
$ws.Range("E15").Value = 5800

# Total Operating Expenses
$ws.Range("D17").Value = 112200
$ws.Range("E17").Value = 114600

# Operating Income or Loss
$ws.Range("D18").Value = 1000
$ws.Range("E18").Value = 2100

# Total Other Income/Expenses Net
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 0

# Earnings Before Interest And Taxes
$ws.Range("D21").Value = 6400
$ws.Range("E21").Value = 7900

# Interest Expense
$ws.Range("D22").Value = 900
$ws.Range("E22").Value = 1100

# Income Before Tax
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 1000

# Income Tax Expense
$ws.Range("D24").Value = 300
$ws.Range("E24").Value = 0

# Minority Interest
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0

# Income After Tax
$ws.Range("D26").Value = -200
$ws.Range("E26").Value = 1100

# Net Income From Continuing Ops
$ws.Range("D27").Value = -200
$ws.Range("E27").Value = 1100

# Non-recurring Events
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0

# Discontinued Operations
$ws.Range("D29").Value = 200
$ws.Range("E29").Value = "NA"

# Extraordinary Items
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0

# Effect Of Accounting Changes
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0

# Other Items
$ws.Range("D32").Value = 0
$ws.Range("E32").Value = 0

# Net Income
$ws.Range("D33").Value = 0
$ws.Range("E33").Value = 1100

# Preferred Stock And Other Adjustments
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0

# Net Income Applicable To Common Shares
$ws.Range("D35").Value = 0
$ws.Range("E35").Value = 1100

# Balance Sheet - Period Ending (row 38)
$ws.Range("D38").Value = 43466
$ws.Range("E38").Value = 43375

# Cash And Cash Equivalents
$ws.Range("D41").Value = 4700
$ws.Range("E41").Value = 1900

# Short Term Investments
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0

# Net Receivables
$ws.Range("D43").Value = 2600
$ws.Range("E43").Value = 2100

# Inventory
$ws.Range("D44").Value = 9600
$ws.Range("E44").Value = 9900

# Other Current Assets
$ws.Range("D45").Value = 6500
$ws.Range("E45").Value = 6700

# Total Current Assets
$ws.Range("D46").Value = 23400
$ws.Range("E46").Value = 20600

# Long Term Investments
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0

# Property Plant and Equipment
$ws.Range("D48").Value = 138800
$ws.Range("E48").Value = 141400

# Goodwill
$ws.Range("D49").Value = 7700
$ws.Range("E49").Value = 7800

# Intangible Assets
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0

# Accumulated Amortization
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0

# Other Assets
$ws.Range("D52").Value = 2200
$ws.Range("E52").Value = 2600

# Deferred Long Term Asset Charges
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0

# Total Assets
$ws.Range("D54").Value = 172000
$ws.Range("E54").Value = 172300

# Accounts Payable
$ws.Range("D57").Value = 7900
$ws.Range("E57").Value = 7300

# Short/Current Long Term Debt
$ws.Range("D58").Value = 700
$ws.Range("E58").Value = 700

# Other Current Liabilities
$ws.Range("D59").Value = 24600
$ws.Range("E59").Value = 22100

# Total Current Liabilities
$ws.Range("D60").Value = 33100
$ws.Range("E60").Value = 30100

# Long Term Debt
$ws.Range("D61").Value = 44200
$ws.Range("E61").Value = 47100

# Other Liabilities
$ws.Range("D62").Value = 42000
$ws.Range("E62").Value = 43100

# Deferred Long Term Liability Charges
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0

# Negative Goodwill (row 64, reuses "Minority Interest" label)
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0

# Total Liabilities
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0

# Total Liabilities value row
$ws.Range("D66").Value = 119400
$ws.Range("E66").Value = 120400

# Misc Stocks Options Warrants
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0

# Redeemable Preferred Stock
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0

# Preferred Stock
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0

# Common Stock
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0

# Retained Earnings
$ws.Range("D72").Value = -111100
$ws.Range("E72").Value = -111200

# Treasury Stock
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0

# Capital Surplus
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0

# Other Stockholder Equity
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0

# Total Stockholder Equity
$ws.Range("D76").Value = 52700
$ws.Range("E76").Value = 52000

# Net Tangible Assets
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0

# Cash Flow Statement - Period Ending (row 80)
$ws.Range("D80").Value = 43466
$ws.Range("E80").Value = 43375

# Net Income (row 81)
$ws.Range("D81").Value = 0
$ws.Range("E81").Value = 1100

# Depreciation
$ws.Range("D83").Value = 5500
$ws.Range("E83").Value = 5800

# Adjustments To Net Income
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0

# Changes In Accounts Receivables
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0

# Changes In Liabilities
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0

# Changes In Inventories
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0

# Changes In Other Operating Activities
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0

# Total Cash Flow From Operating Activities
$ws.Range("D89").Value = 10200
$ws.Range("E89").Value = -7700

# Capital Expenditures
$ws.Range("D91").Value = -4400
$ws.Range("E91").Value = -2900

# Investments
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0

# Other Cashflows from Investing Activities
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0

# Total Cash Flows From Investing Activities
$ws.Range("D94").Value = -4400
$ws.Range("E94").Value = -2400

# Dividends Paid
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0

# Sale Purchase of Stock
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0

# Net Borrowings
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0

# Other Cash Flows from Financing Activities
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0

# Total Cash Flows From Financing Activities
$ws.Range("D100").Value = -3100
$ws.Range("E100").Value = 8400

# Effect Of Exchange Rate Changes
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0

# Change In Cash and Cash Equivalents
$ws.Range("D102").Value = 2700
$ws.Range("E102").Value = -1700
